# Bug fix for one trial with wrong date
# Updates computed statistics (publications within 24m, percentage, and
# confidence interval bounds) for the affected institution rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Aalborg University Hospital
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 54.3
$ws.Range("E3").Value = 38.2
$ws.Range("F3").Value = 69.5

# Row 4 - Aarhus University
$ws.Range("C4").Value = 113
$ws.Range("D4").Value = 60.8
$ws.Range("E4").Value = 53.6
$ws.Range("F4").Value = 67.5

# Row 44 - University of Copenhagen
$ws.Range("C44").Value = 43
$ws.Range("D44").Value = 44.3
$ws.Range("E44").Value = 34.8
$ws.Range("F44").Value = 54.2

# Row 51 - University of Tampere
$ws.Range("C51").Value = 6
$ws.Range("D51").Value = 66.7
$ws.Range("E51").Value = 35.4
$ws.Range("F51").Value = 87.9
